$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add "Wins", "Losses", "Ties" in AD1:AF1, matching the
#     bold/centered/bordered style already used by the other header cells.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the formatting (not the value) from an existing header cell (A1,
# style index 1) onto the three new header cells so they reuse the same
# style instead of creating new ones.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# --- Data rows 2-50: season record columns (Wins=80, Losses=82, Ties=0)
for ($r = 2; $r -le 50; $r++) {
    $ws.Cells.Item($r, 30).Value = 80
    $ws.Cells.Item($r, 31).Value = 82
    $ws.Cells.Item($r, 32).Value = 0
}
